$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 4562.16
$ws.Range("I40").Value = 4135.2666
$ws.Range("J40").Value = 5202.5
$ws.Range("K40").Value = 4135.2666
$ws.Range("L40").Value = 5202.5
$ws.Range("M40").Value = -3960.2666
$ws.Range("N40").Value = -5552.5
# Row 62
$ws.Range("H62").Value = 10766.777
$ws.Range("I62").Value = 6912
$ws.Range("J62").Value = 11248.625
$ws.Range("K62").Value = 6912
$ws.Range("L62").Value = 11248.625
$ws.Range("M62").Value = -6288
$ws.Range("N62").Value = -12496.625
# Row 65
$ws.Range("H65").Value = 10766.777
$ws.Range("I65").Value = 6912
$ws.Range("J65").Value = 11248.625
$ws.Range("K65").Value = 34560
$ws.Range("L65").Value = 56243.125
$ws.Range("M65").Value = -31440
$ws.Range("N65").Value = -62483.125
# Row 116
$ws.Range("H116").Value = 7958.8
$ws.Range("I116").Value = 10364.667
$ws.Range("J116").Value = 4350
$ws.Range("K116").Value = 10364.667
$ws.Range("L116").Value = 4350
$ws.Range("M116").Value = -6922.666999999999
$ws.Range("N116").Value = -11234
# Row 138
$ws.Range("H138").Value = 2438.5
$ws.Range("I138").Value = 1742.7
$ws.Range("J138").Value = 3134.3
$ws.Range("K138").Value = 5228.1
$ws.Range("L138").Value = 9402.900000000001
$ws.Range("M138").Value = -88.10000000000036
$ws.Range("N138").Value = -19682.9

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2278006
$ws.Range("I32").Value = 3380.5789
$ws.Range("J32").Value = 16683966
$ws.Range("K32").Value = 3380.5789
$ws.Range("L32").Value = 16683966
$ws.Range("M32").Value = -3093.5789
$ws.Range("N32").Value = -16684540
# Row 61
$ws.Range("H61").Value = 4501.75
$ws.Range("I61").Value = 4166.6665
$ws.Range("J61").Value = 5507
$ws.Range("K61").Value = 4166.6665
$ws.Range("L61").Value = 5507
$ws.Range("M61").Value = -3954.6665
$ws.Range("N61").Value = -5931
# Row 74
$ws.Range("I74").Value = 7500
$ws.Range("J74").Value = 9000
$ws.Range("K74").Value = 7500
$ws.Range("L74").Value = 9000
$ws.Range("M74").Value = -6626
$ws.Range("N74").Value = -10748
# Row 77
$ws.Range("I77").Value = 7500
$ws.Range("J77").Value = 9000
$ws.Range("K77").Value = 37500
$ws.Range("L77").Value = 45000
$ws.Range("M77").Value = -33132
$ws.Range("N77").Value = -53736
# Row 136
$ws.Range("H136").Value = 4501.75
$ws.Range("I136").Value = 4166.6665
$ws.Range("J136").Value = 5507
$ws.Range("K136").Value = 12499.9995
$ws.Range("L136").Value = 16521
$ws.Range("M136").Value = -9949.999500000002
$ws.Range("N136").Value = -21621

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2363.476
$ws.Range("I134").Value = 2035.2142
$ws.Range("J134").Value = 3020
$ws.Range("K134").Value = 6105.642599999999
$ws.Range("L134").Value = 9060
$ws.Range("M134").Value = -3570.642599999999
$ws.Range("N134").Value = -14130

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 5595.778
$ws.Range("I58").Value = 999
$ws.Range("J58").Value = 6170.375
$ws.Range("K58").Value = 999
$ws.Range("L58").Value = 6170.375
$ws.Range("M58").Value = -796
$ws.Range("N58").Value = -6576.375
# Row 106
$ws.Range("H106").Value = 26720.166
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 26720.166
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 26720.166
$ws.Range("N106").Value = -29244.166
# Row 136
$ws.Range("H136").Value = 5595.778
$ws.Range("I136").Value = 999
$ws.Range("J136").Value = 6170.375
$ws.Range("K136").Value = 2997
$ws.Range("L136").Value = 18511.125
$ws.Range("M136").Value = -447
$ws.Range("N136").Value = -23611.125

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 4200
$ws.Range("I5").Value = 4200
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 12600
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -12488
$ws.Range("N5").Value = ""
# Row 14
$ws.Range("H14").Value = 332.83334
$ws.Range("I14").Value = 332.83334
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 998.5000200000001
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -825.5000200000001
# Row 109
$ws.Range("H109").Value = 1003070
$ws.Range("I109").Value = 3334333.2
$ws.Range("J109").Value = 3957.1428
$ws.Range("K109").Value = 10002999.6
$ws.Range("L109").Value = 11871.4284
$ws.Range("M109").Value = -10001959.6
$ws.Range("N109").Value = -13951.4284
# Row 135
$ws.Range("H135").Value = 4200
$ws.Range("I135").Value = 4200
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 37800
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -35265
$ws.Range("N135").Value = ""

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 3496.6667
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 4495
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 13485
$ws.Range("M126").Value = -2030
$ws.Range("N126").Value = -18425
# Row 132
$ws.Range("H132").Value = 4281.7075
$ws.Range("I132").Value = 3778.5151
$ws.Range("J132").Value = 6357.375
$ws.Range("K132").Value = 11335.5453
$ws.Range("L132").Value = 19072.125
$ws.Range("M132").Value = -8805.5453
$ws.Range("N132").Value = -24132.125
# Row 134
$ws.Range("H134").Value = 100274.375
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 100274.375
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 300823.125
$ws.Range("N134").Value = -305893.125

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 695
$ws.Range("I16").Value = 695
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 695
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -525
# Row 69
$ws.Range("H69").Value = 33331
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 33331
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 33331
$ws.Range("N69").Value = -34953
# Row 72
$ws.Range("H72").Value = 33331
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 33331
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 99993
$ws.Range("N72").Value = -108105
# Row 93
$ws.Range("H93").Value = 1086.8572
$ws.Range("I93").Value = 1234.6666
$ws.Range("J93").Value = 200
$ws.Range("K93").Value = 1234.6666
$ws.Range("L93").Value = 200
$ws.Range("M93").Value = 13.33339999999998
$ws.Range("N93").Value = -2696
# Row 132
$ws.Range("H132").Value = 2969.7273
$ws.Range("I132").Value = 2916.7
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 8750.099999999999
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -6220.099999999999
$ws.Range("N132").Value = -15560
# Row 136
$ws.Range("H136").Value = 4114
$ws.Range("I136").Value = 2837
$ws.Range("J136").Value = 6029.5
$ws.Range("K136").Value = 8511
$ws.Range("L136").Value = 18088.5
$ws.Range("M136").Value = -5961
$ws.Range("N136").Value = -23188.5
# Row 137
$ws.Range("H137").Value = 750000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 750000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 750000
$ws.Range("M137").Value = ""
$ws.Range("N137").Value = -760200

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 994.25
$ws.Range("I113").Value = 1290.8
$ws.Range("J113").Value = 500
$ws.Range("K113").Value = 3872.4
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = -1702.4
$ws.Range("N113").Value = -5840
# Row 126
$ws.Range("H126").Value = 7250
$ws.Range("I126").Value = 3500
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 10500
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -8030
$ws.Range("N126").Value = -28940
# Row 136
$ws.Range("H136").Value = 9994.5
$ws.Range("I136").Value = 9994
$ws.Range("J136").Value = 9995
$ws.Range("K136").Value = 29982
$ws.Range("L136").Value = 29985
$ws.Range("M136").Value = -27432
$ws.Range("N136").Value = -35085
